$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.383.25'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '2.893.02'
$ws.Range("E3").Value = '  -3.88%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.43'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.503'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("D9").Value = '2.891.49'
$ws.Range("E9").Value = '  -3.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.14%  '
$ws.Range("E11").Value = '  -3.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.07'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '3.371.75'
$ws.Range("E16").Value = '  -4.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.79'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D18").Value = '60.377.12'
$ws.Range("E18").Value = '  -2.81%  '
$ws.Range("D19").Value = '2.892.29'
$ws.Range("E19").Value = '  -3.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '424.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.668'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.10'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.83'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.91'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("E26").Value = '  -3.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.74'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.26'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.17'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.43'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.55%  '
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.65'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.21'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.73'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.20%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.122'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.286'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.05'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0344'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '370.23'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.24%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.09'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.646.81'
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.93'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.01%  '
$ws.Range("E51").Value = '  -1.15%  '
